$d = $word.ActiveDocument

# 1. Insert a new paragraph ("asdfasdf") before the very first paragraph
#    ("Conditional Display"). InsertParagraphBefore() clones the
#    paragraph formatting (pStyle "8" + spacing) of the following
#    paragraph automatically.
$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphBefore()
$newFirstPara = $d.Paragraphs(1)
$newFirstPara.Range.Text = "asdfasdf"

# 2. Add a run of text ("asdfasdf") right before the bookmarkStart/
#    bookmarkEnd in the final (empty) paragraph of the document.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertBefore("asdfasdf")

# 3. Add top/bottom cell margins (0 dxa) to the table's tblCellMar,
#    next to the existing left/right margins.
$table = $d.Tables(1)
$table.TopPadding = 0
$table.BottomPadding = 0

# 4. Mark the "Table Grid" table style as a Quick Style (adds
#    <w:qFormat/> to its style definition).
foreach ($s in $d.Styles) {
    if ($s.NameLocal -eq "Table Grid") {
        $s.QuickStyle = $true
    }
}

Write-Output "done"
